$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.0001240510425049024
$ws.Range("C2").Value = -0.00001133570006395246

$ws.Range("B3").Value = -0.00563604457062183428
$ws.Range("C3").Value = -0.00162538343594320601

$ws.Range("C4").Value = -0.04135662873497381042

$ws.Range("B5").Value = -0.00021305533987048259
$ws.Range("C5").Value = -0.00001946885234516138

$ws.Range("B6").Value = -0.0001888587432858913
$ws.Range("C6").Value = -0.00001725778095362784

$ws.Range("B7").Value = -0.07294965438723011175
$ws.Range("C7").Value = -0.00666357887251933789
